$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 395, shifting existing rows 395:413 down to 396:414
$ws.Rows.Item(395).Insert()

# Populate the newly inserted row 395 with the new record
$ws.Cells.Item(395, 1).Value = 4
$ws.Cells.Item(395, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(395, 3).Value = "Los Lagos"
$ws.Cells.Item(395, 4).Value = 44939
$ws.Cells.Item(395, 5).Value = 10
$ws.Cells.Item(395, 6).Value = 100112045
$ws.Cells.Item(395, 7).Value = "Zapallo"
$ws.Cells.Item(395, 8).Value = "Paine"
$ws.Cells.Item(395, 9).Value = "1a (cosecha)"
$ws.Cells.Item(395, 10).Value = 1200
$ws.Cells.Item(395, 11).Value = 650
$ws.Cells.Item(395, 12).Value = 650
$ws.Cells.Item(395, 13).Value = 650
$ws.Cells.Item(395, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(395, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(395, 16).Value = 650
$ws.Cells.Item(395, 17).Value = 1
$ws.Cells.Item(395, 18).Value = "Hortaliza"
